# Adjusted risk calc formula
#
# The "LF70" watershed-risk row's Total Risk / Future Risk values change,
# which in turn bumps its "LF3" neighbour's Total Risk / Current Risk /
# Future Risk high enough that it swaps rank order with the "LF1" row
# below it (their LF description + rank/current/future values trade
# places), and the following "LF39" row's Rank shifts down by one to
# stay consistent.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row for "LF70: ... small population size ..." (Word table row 13)
$t.Cell(13, 4).Range.Text = "6"
$t.Cell(13, 6).Range.Text = "M"

# Row that was "LF1: ... predation from pinnipeds ..." (Word table row 14)
# becomes "LF3: ... anthropogenic activity (non fishing) ..."
$t.Cell(14, 2).Range.Text = "LF3: Mortality or fitness reduction as a result of stress due to anthropogenic activity (non fishing)"
$t.Cell(14, 4).Range.Text = "3"
$t.Cell(14, 5).Range.Text = "VL"
$t.Cell(14, 6).Range.Text = "M"

# Row that was "LF3: ... anthropogenic activity (non fishing) ..." (Word
# table row 15) becomes "LF1: ... predation from pinnipeds ..."
$t.Cell(15, 2).Range.Text = "LF1: Mortality or fitness reduction due to predation from pinnipeds or other aquatic species"
$t.Cell(15, 3).Range.Text = "14"
$t.Cell(15, 5).Range.Text = "L"
$t.Cell(15, 6).Range.Text = "VL"

# Row for "LF39: ... stranding in rearing habitat ..." (Word table row 16)
$t.Cell(16, 3).Range.Text = "14"
